$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the bold/bordered header
# style already used by the other header cells (e.g. G1 "sum").
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"
$excel.CutCopyMode = $false

# Fill in the Save values for each existing data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
